# Update "Pais" COVID data sheet: refresh case counts and re-rank countries
# whose totals changed enough to swap places in the (descending by
# "Casos totales") ranking, plus bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- timestamp (row 1) -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 16:16"

# --- straightforward numeric refreshes (no re-ranking needed) --------
# Estados Unidos (row 6)
$ws.Range("B6").Value = 46285
$ws.Range("C6").Value = 2551
$ws.Range("E6").Value = 45402
$ws.Range("G6").Value = 35
$ws.Range("H6").Value = 588

# Austria (row 15)
$ws.Range("B15").Value = 4971
$ws.Range("C15").Value = 497
$ws.Range("E15").Value = 4934

# Noruega (row 17)
$ws.Range("B17").Value = 2753
$ws.Range("C17").Value = 128
$ws.Range("E17").Value = 2735

# --- Chile / Pakistan swap places (rows 31-32) ------------------------
# Pakistan's new total (958) overtakes Chile's (922), so Pakistan moves
# up to row 31 and Chile drops to row 32.
$ws.Range("A31").Value = "Pakistan"
$ws.Range("B31").Value = 958
$ws.Range("C31").Value = 83
$ws.Range("D31").Value = 13
$ws.Range("E31").Value = 938
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 7

$ws.Range("A32").Value = "Chile"
$ws.Range("B32").Value = 922
$ws.Range("C32").Value = 176
$ws.Range("D32").Value = 17
$ws.Range("E32").Value = 903
$ws.Range("F32").Value = 7
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 2

# --- Ucrania jumps ahead of Camboya / Azerbaiyan / Senegal (rows 93-96)
# Ucrania's new total (97) pushes it above the other three, which each
# shift down one row.
$ws.Range("A93").Value = "Ucrania"
$ws.Range("B93").Value = 97
$ws.Range("C93").Value = 24
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = 93
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 3

$ws.Range("A94").Value = "Camboya"
$ws.Range("B94").Value = 91
$ws.Range("C94").Value = 4
$ws.Range("D94").Value = 4
$ws.Range("E94").Value = 87
$ws.Range("F94").Value = 1
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 0

$ws.Range("A95").Value = "Azerbaiyan"
$ws.Range("B95").Value = 87
$ws.Range("C95").Value = 15
$ws.Range("D95").Value = 10
$ws.Range("E95").Value = 76
$ws.Range("F95").Value = 6
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 1

$ws.Range("A96").Value = "Senegal"
$ws.Range("B96").Value = 86
$ws.Range("C96").Value = 7
$ws.Range("D96").Value = 8
$ws.Range("E96").Value = 78
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0

# --- Macao jumps ahead of Kenia / Costa de Marfil (rows 126-128) ------
$ws.Range("A126").Value = "Macao"
$ws.Range("B126").Value = 26
$ws.Range("C126").Value = 1
$ws.Range("D126").Value = 10
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 0

$ws.Range("A127").Value = "Kenia"
$ws.Range("B127").Value = 25
$ws.Range("C127").Value = 9
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 25
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0

$ws.Range("A128").Value = "Costa de Marfil"
$ws.Range("B128").Value = 25
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 2
$ws.Range("E128").Value = 23
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0
